$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Basic Version Of Profile Screen."
#   - adds a new "Post Views" service row (WS-PS-10) to the Post section
#   - adds a new "Profile" service row (WS-UP-10) to the User section
# Two blank rows are inserted first (top to bottom) so the remaining rows
# settle in their final positions; the new cell values are written
# afterwards (Post Views before Profile) so new shared-string entries come
# out in the same order as the authoritative edit.
# ---------------------------------------------------------------------------

# --- structural inserts -----------------------------------------------------
$ws.Range("B15").EntireRow.Insert()
$ws.Range("B16:O16").Copy()
$ws.Range("B15:O15").PasteSpecial(-4122) | Out-Null

# Re-apply the filter/defined-name range here (B3:L32 - i.e. +1 row for the
# Profile insert above) *before* the second row insert below, so that the
# Post Views insert naturally grows the filter range by one more row instead
# of the AutoFilter call snapping to the full current region afterwards.
$ws.AutoFilterMode = $false
$ws.Range("B3:L32").AutoFilter() | Out-Null

$ws.Range("B31").EntireRow.Insert()
$ws.Range("B32:O32").Copy()
$ws.Range("B31:O31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- row 31: new "Post Views" service (Post section) -----------------------
$ws.Range("D31").Value2 = "WS-PS-10"
$ws.Range("H31").Value2 = "/views"
$ws.Range("E31").Value2 = "app.post.save.viewe"
$ws.Range("C31").Value2 = "Post Views"
$ws.Range("B31").Value2 = "Post"
$ws.Range("F31").Value2 = $false
$ws.Range("G31").Value2 = "post"
$ws.Range("I31").Value2 = "POST"
$ws.Range("K31").Value2 = "Done"
$ws.Range("M31").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D31,"'',''CONNON_CONFIG'', 0, ''",C31,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N31").Formula = '=_xlfn.CONCAT(IF(I31="GET","@GetMapping(",IF(I31="POST","@PostMapping(",IF(I31="DELETE","@DeleteMapping(",IF(I31="PUT","@PutMapping(","")))),CHAR(34),H31,CHAR(34),")")'
$ws.Range("O31").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D31,,CHAR(34),", serviceName = ",CHAR(34),C31,CHAR(34), ", queryId = ",CHAR(34),E31,CHAR(34),", logActivity =",F31,")")'

# --- row 15: new "Profile" service (User section) --------------------------
$ws.Range("C15").Value2 = "Profile"
$ws.Range("D15").Value2 = "WS-UP-10"
$ws.Range("E15").Value2 = "app.user.profile.get"
$ws.Range("H15").Value2 = "/profle"
$ws.Range("B15").Value2 = "User"
$ws.Range("F15").Value2 = $false
$ws.Range("G15").Value2 = "user"
$ws.Range("I15").Value2 = "POST"
$ws.Range("J15").Value2 = "P1"
$ws.Range("M15").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D15,"'',''CONNON_CONFIG'', 0, ''",C15,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N15").Formula = '=_xlfn.CONCAT(IF(I15="GET","@GetMapping(",IF(I15="POST","@PostMapping(",IF(I15="DELETE","@DeleteMapping(",IF(I15="PUT","@PutMapping(","")))),CHAR(34),H15,CHAR(34),")")'
$ws.Range("O15").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D15,,CHAR(34),", serviceName = ",CHAR(34),C15,CHAR(34), ", queryId = ",CHAR(34),E15,CHAR(34),", logActivity =",F15,")")'

# --- defined name mirrors the (now B3:L32) filter range --------------------
$fdb = $wb.Names.Item("ServicesList!_FilterDatabase")
$fdb.RefersTo = "=ServicesList!`$B`$3:`$L`$32"

# --- view state: scrolled down a bit, selection on the new Profile row -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("J15").Select() | Out-Null
